$wb = $excel.ActiveWorkbook

$wsEnroll = $wb.Worksheets.Item("enrollment_dd")
$wsFollow = $wb.Worksheets.Item("followup_dd")
$wsSwf    = $wb.Worksheets.Item("swf_dd")

# --- enrollment_dd: pregnant_date (row 8) gets a LowerRange / UpperRange date window ---
# Set H8 before G8 so the new shared strings are appended in "+0d","-9m" order.
$wsEnroll.Range("H8").Value = "'+0d"
$wsEnroll.Range("G8").Value = "'-9m"

# enrollment_dd row 12 (device_use): allow Don't know / Refuse
$wsEnroll.Range("J12").Value = $true
$wsEnroll.Range("K12").Value = $true

# --- followup_dd: vdate (row 4) gets a LowerRange / UpperRange date window ---
$wsFollow.Range("G4").Value = "'-2m"
$wsFollow.Range("H4").Value = "'+0d"

# --- swf_dd: withdrawdate (row 8) gets a LowerRange / UpperRange date window ---
# First bring G8/H8's formatting in line with the other date-range cells (copy from followup_dd!G4),
# then fill in the values (quote-prefixed, since they start with +/-).
$wsFollow.Range("G4").Copy()
$wsSwf.Range("G8").PasteSpecial(-4122)
$wsSwf.Range("H8").PasteSpecial(-4122)
$wsSwf.Range("G8").Value = "'-1y"
$wsSwf.Range("H8").Value = "'+0d"
$wsSwf.Range("J8").Value = $true

# Row 11 on swf_dd is no longer using the maximum auto height; it now uses a fixed custom height.
$wsSwf.Rows.Item(11).RowHeight = 184.2

# --- Restore/update the active-cell selections on each sheet, matching where the author left off ---
$wsEnroll.Activate()
$wsEnroll.Range("J13").Select() | Out-Null

$wsFollow.Activate()
$wsFollow.Range("G11").Select() | Out-Null

$wsSwf.Activate()
$wsSwf.Range("J9").Select() | Out-Null

$wsEnroll.Activate()
